$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 2145.6667
$ws.Range("I42").Value = 400
$ws.Range("J42").Value = 2494.8
$ws.Range("K42").Value = 1200
$ws.Range("L42").Value = 7484.400000000001
$ws.Range("M42").Value = -970
$ws.Range("N42").Value = -7944.400000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1437
$ws.Range("I62").Value = 1600
$ws.Range("J62").Value = 1111
$ws.Range("K62").Value = 1600
$ws.Range("L62").Value = 1111
$ws.Range("M62").Value = -976
$ws.Range("N62").Value = -2359

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1437
$ws.Range("I65").Value = 1600
$ws.Range("J65").Value = 1111
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 5555
$ws.Range("M65").Value = -4880
$ws.Range("N65").Value = -11795

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3619.9473
$ws.Range("I98").Value = 3269.7
$ws.Range("J98").Value = 4009.111
$ws.Range("K98").Value = 3269.7
$ws.Range("L98").Value = 4009.111
$ws.Range("M98").Value = -1771.7
$ws.Range("N98").Value = -7005.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 9271.111000000001
$ws.Range("I106").Value = 9271.111000000001
$ws.Range("K106").Value = 9271.111000000001
$ws.Range("M106").Value = -8640.111000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1474.7693
$ws.Range("I107").Value = 1222.3
$ws.Range("K107").Value = 1222.3
$ws.Range("M107").Value = 697.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4403.8335
$ws.Range("J112").Value = 6059.125
$ws.Range("L112").Value = 18177.375
$ws.Range("N112").Value = -20393.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3619.9473
$ws.Range("I122").Value = 3269.7
$ws.Range("J122").Value = 4009.111
$ws.Range("K122").Value = 9809.099999999999
$ws.Range("L122").Value = 12027.333
$ws.Range("M122").Value = -7359.099999999999
$ws.Range("N122").Value = -16927.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 7086.909
$ws.Range("I135").Value = 1682.1875
$ws.Range("K135").Value = 15139.6875
$ws.Range("M135").Value = -12604.6875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2565.0461
$ws.Range("I138").Value = 1421.12
$ws.Range("J138").Value = 3280
$ws.Range("K138").Value = 4263.36
$ws.Range("L138").Value = 9840
$ws.Range("M138").Value = 876.6400000000003
$ws.Range("N138").Value = -20120

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2930.125
$ws.Range("I141").Value = 2854.4092
$ws.Range("J141").Value = 3763
$ws.Range("K141").Value = 8563.2276
$ws.Range("L141").Value = 11289
$ws.Range("M141").Value = -3383.2276
$ws.Range("N141").Value = -21649

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2035.84
$ws.Range("I32").Value = 908.86566
$ws.Range("K32").Value = 908.86566
$ws.Range("M32").Value = -621.86566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2420.318
$ws.Range("I45").Value = 2183.8235
$ws.Range("K45").Value = 2183.8235
$ws.Range("M45").Value = -1806.8235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3465.2
$ws.Range("J88").Value = 3465.2
$ws.Range("L88").Value = 3465.2
$ws.Range("N88").Value = -4277.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3465.2
$ws.Range("J91").Value = 3465.2
$ws.Range("L91").Value = 3465.2
$ws.Range("N91").Value = -6273.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3520.9524
$ws.Range("I122").Value = 3291.2354
$ws.Range("J122").Value = 4497.25
$ws.Range("K122").Value = 9873.706200000001
$ws.Range("L122").Value = 13491.75
$ws.Range("M122").Value = -7423.706200000001
$ws.Range("N122").Value = -18391.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12504478
$ws.Range("J86").Value = 5485.222
$ws.Range("L86").Value = 5485.222
$ws.Range("N86").Value = -7731.222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 12504478
$ws.Range("J89").Value = 5485.222
$ws.Range("L89").Value = 27426.11
$ws.Range("N89").Value = -38658.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2707.9473
$ws.Range("I99").Value = 2075.6428
$ws.Range("K99").Value = 2075.6428
$ws.Range("M99").Value = -577.6428000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14067
$ws.Range("I16").Value = 1082.1111
$ws.Range("J16").Value = 72499
$ws.Range("K16").Value = 1082.1111
$ws.Range("L16").Value = 72499
$ws.Range("M16").Value = -795.1111000000001
$ws.Range("N16").Value = -73073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3275.5254
$ws.Range("I31").Value = 5081.8423
$ws.Range("K31").Value = 5081.8423
$ws.Range("M31").Value = -4786.8423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3275.5254
$ws.Range("I34").Value = 5081.8423
$ws.Range("K34").Value = 5081.8423
$ws.Range("M34").Value = -4879.8423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 37640.2
$ws.Range("J92").Value = 42675.25
$ws.Range("L92").Value = 42675.25
$ws.Range("N92").Value = -47667.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9265.416999999999
$ws.Range("I99").Value = 5338.8
$ws.Range("J99").Value = 12070.143
$ws.Range("K99").Value = 5338.8
$ws.Range("L99").Value = 12070.143
$ws.Range("M99").Value = -3840.8
$ws.Range("N99").Value = -15066.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 14067
$ws.Range("I113").Value = 1082.1111
$ws.Range("J113").Value = 72499
$ws.Range("K113").Value = 1082.1111
$ws.Range("L113").Value = 72499
$ws.Range("M113").Value = 1087.8889
$ws.Range("N113").Value = -76839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1175.6111
$ws.Range("I122").Value = 1014.4167
$ws.Range("K122").Value = 3043.2501
$ws.Range("M122").Value = -593.2501000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9265.416999999999
$ws.Range("I126").Value = 5338.8
$ws.Range("J126").Value = 12070.143
$ws.Range("K126").Value = 16016.4
$ws.Range("L126").Value = 36210.429
$ws.Range("M126").Value = -13546.4
$ws.Range("N126").Value = -41150.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 71995.60000000001
$ws.Range("J131").Value = 71995.60000000001
$ws.Range("L131").Value = 71995.60000000001
$ws.Range("N131").Value = -82075.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 705.05554
$ws.Range("I134").Value = 678.3939
$ws.Range("K134").Value = 2035.1817
$ws.Range("M134").Value = 499.8182999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69949
$ws.Range("J37").Value = 69949
$ws.Range("L37").Value = 209847
$ws.Range("N37").Value = -210071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2040
$ws.Range("J41").Value = 5000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7769.087
$ws.Range("I56").Value = 7769.087
$ws.Range("K56").Value = 7769.087
$ws.Range("M56").Value = -7239.087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2779.4
$ws.Range("I86").Value = 2779.4
$ws.Range("K86").Value = 8338.200000000001
$ws.Range("M86").Value = -7152.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 2779.4
$ws.Range("I89").Value = 2779.4
$ws.Range("K89").Value = 25014.6
$ws.Range("M89").Value = -19086.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 899.2
$ws.Range("I97").Value = 999
$ws.Range("K97").Value = 2997
$ws.Range("M97").Value = -2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4284.143
$ws.Range("I134").Value = 1664.8334
$ws.Range("K134").Value = 4994.5002
$ws.Range("M134").Value = 75.4997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8397983
$ws.Range("I11").Value = 5501446.5
$ws.Range("K11").Value = 5501446.5
$ws.Range("M11").Value = -5501307.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 765.1111
$ws.Range("I97").Value = 727
$ws.Range("J97").Value = 898.5
$ws.Range("K97").Value = 727
$ws.Range("L97").Value = 898.5
$ws.Range("M97").Value = -231
$ws.Range("N97").Value = -1890.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 23254.889
$ws.Range("I113").Value = 25161.75
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 25161.75
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -22991.75
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3333.05
$ws.Range("I126").Value = 2734.5
$ws.Range("J126").Value = 4729.6665
$ws.Range("K126").Value = 8203.5
$ws.Range("L126").Value = 14188.9995
$ws.Range("M126").Value = -5733.5
$ws.Range("N126").Value = -19128.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3274.6155
$ws.Range("I61").Value = 3007.7778
$ws.Range("J61").Value = 3875
$ws.Range("K61").Value = 3007.7778
$ws.Range("L61").Value = 3875
$ws.Range("M61").Value = -2805.7778
$ws.Range("N61").Value = -4279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3274.6155
$ws.Range("I113").Value = 3007.7778
$ws.Range("J113").Value = 3875
$ws.Range("K113").Value = 3007.7778
$ws.Range("L113").Value = 3875
$ws.Range("M113").Value = -837.7777999999998
$ws.Range("N113").Value = -8215

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9872.277
$ws.Range("I122").Value = 11499.214
$ws.Range("K122").Value = 34497.642
$ws.Range("M122").Value = -32047.642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2817.7222
$ws.Range("I132").Value = 2754.7144
$ws.Range("K132").Value = 8264.143199999999
$ws.Range("M132").Value = -5734.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12270.454
$ws.Range("I62").Value = 12397.7
$ws.Range("K62").Value = 12397.7
$ws.Range("M62").Value = -11773.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 12270.454
$ws.Range("I65").Value = 12397.7
$ws.Range("K65").Value = 61988.5
$ws.Range("M65").Value = -58868.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3463
$ws.Range("I81").Value = 3463
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6926
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5865
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3463
$ws.Range("I84").Value = 3463
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 34630
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -29326
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1588.0286
$ws.Range("I136").Value = 1457.5161
$ws.Range("K136").Value = 4372.5483
$ws.Range("M136").Value = -1822.5483
